$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 517, shifting all subsequent
# rows (517-541) down to (519-543).
$ws.Rows("517:518").Insert()

# New row 517 data
$ws.Cells.Item(517,1).Value = 6
$ws.Cells.Item(517,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(517,3).Value = "Metropolitana"
$ws.Cells.Item(517,4).Value = 44516
$ws.Cells.Item(517,5).Value = 13
$ws.Cells.Item(517,6).Value = 100112031
$ws.Cells.Item(517,7).Value = "Poroto verde"
$ws.Cells.Item(517,8).Value = "Magnum"
$ws.Cells.Item(517,9).Value = "Primera"
$ws.Cells.Item(517,10).Value = 120
$ws.Cells.Item(517,11).Value = 45000
$ws.Cells.Item(517,12).Value = 47000
$ws.Cells.Item(517,13).Value = 46167
$ws.Cells.Item(517,14).Value = "$/caja 25 kilos"
$ws.Cells.Item(517,15).Value = "Provincia de Limarí"
$ws.Cells.Item(517,16).Value = 1847
$ws.Cells.Item(517,17).Value = 25
$ws.Cells.Item(517,18).Value = "Hortaliza"

# New row 518 data
$ws.Cells.Item(518,1).Value = 6
$ws.Cells.Item(518,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(518,3).Value = "Metropolitana"
$ws.Cells.Item(518,4).Value = 44516
$ws.Cells.Item(518,5).Value = 13
$ws.Cells.Item(518,6).Value = 100112031
$ws.Cells.Item(518,7).Value = "Poroto verde"
$ws.Cells.Item(518,8).Value = "Magnum"
$ws.Cells.Item(518,9).Value = "Primera"
$ws.Cells.Item(518,10).Value = 250
$ws.Cells.Item(518,11).Value = 50000
$ws.Cells.Item(518,12).Value = 55000
$ws.Cells.Item(518,13).Value = 52000
$ws.Cells.Item(518,14).Value = "$/saco 25 kilos"
$ws.Cells.Item(518,15).Value = "Región de O'Higgins"
$ws.Cells.Item(518,16).Value = 2080
$ws.Cells.Item(518,17).Value = 25
$ws.Cells.Item(518,18).Value = "Hortaliza"
